$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @("FAPs", "Efna5", "Epha7", "ECs", 3, 1, 2.030023666666667, 6.090071, 0.8776223887075381, 0.8776223887075382, 2, 0.6666666666666666, 0.03932533333333333, 0.117976, 0.03588667238536742, 0.03588667238536743, 0.07983135736622221, 0.718482216296, 0.031494947141611, 0.03149494714161101)
for ($c = 1; $c -le $row2.Length; $c++) {
    $ws.Cells.Item(2, $c).Value = $row2[$c-1]
}

$row3 = @("FAPs", "Efna5", "Epha7", "FAPs", 3, 1, 2.030023666666667, 6.090071, 0.8776223887075381, 0.8776223887075382, 2, 0.6666666666666666, 0.242411, 0.727233, 0.2212142505155957, 0.2212142505155957, 0.4921000670603334, 4.428900603543, 0.1941425789536448, 0.1941425789536448)
for ($c = 1; $c -le $row3.Length; $c++) {
    $ws.Cells.Item(3, $c).Value = $row3[$c-1]
}

$row4 = @("FAPs", "Efna5", "Epha7", "Neutro", 3, 1, 2.030023666666667, 6.090071, 0.8776223887075381, 0.8776223887075382, 1, 0.3333333333333333, 0.07026666666666666, 0.2108, 0.06412245320095149, 0.06412245320095149, 0.1426429963111111, 1.2837869668, 0.05627530054800636, 0.05627530054800637)
for ($c = 1; $c -le $row4.Length; $c++) {
    $ws.Cells.Item(4, $c).Value = $row4[$c-1]
}

$row5 = @("FAPs", "Efna5", "Epha7", "sCs", 3, 1, 2.030023666666667, 6.090071, 0.8776223887075381, 0.8776223887075382, 3, 1, 0.743817, 2.231451, 0.6787766238980854, 0.6787766238980855, 1.509966113669, 13.589695023021, 0.5957095620642758, 0.5957095620642761)
for ($c = 1; $c -le $row5.Length; $c++) {
    $ws.Cells.Item(5, $c).Value = $row5[$c-1]
}

$row6 = @("sCs", "Efna5", "Epha7", "ECs", 2, 0.6666666666666666, 0.283071, 0.849213, 0.1223776112924619, 0.1223776112924619, 2, 0.6666666666666666, 0.03932533333333333, 0.117976, 0.03588667238536742, 0.03588667238536743, 0.011131861432, 0.100186752888, 0.00439172524375642, 0.00439172524375642)
for ($c = 1; $c -le $row6.Length; $c++) {
    $ws.Cells.Item(6, $c).Value = $row6[$c-1]
}

$row7 = @("sCs", "Efna5", "Epha7", "FAPs", 2, 0.6666666666666666, 0.283071, 0.849213, 0.1223776112924619, 0.1223776112924619, 2, 0.6666666666666666, 0.242411, 0.727233, 0.2212142505155957, 0.2212142505155957, 0.06861952418100001, 0.617575717629, 0.02707167156195085, 0.02707167156195085)
for ($c = 1; $c -le $row7.Length; $c++) {
    $ws.Cells.Item(7, $c).Value = $row7[$c-1]
}

$row8 = @("sCs", "Efna5", "Epha7", "Neutro", 2, 0.6666666666666666, 0.283071, 0.849213, 0.1223776112924619, 0.1223776112924619, 1, 0.3333333333333333, 0.07026666666666666, 0.2108, 0.06412245320095149, 0.06412245320095149, 0.0198904556, 0.1790141004, 0.007847152652945119, 0.007847152652945119)
for ($c = 1; $c -le $row8.Length; $c++) {
    $ws.Cells.Item(8, $c).Value = $row8[$c-1]
}

$row9 = @("sCs", "Efna5", "Epha7", "sCs", 2, 0.6666666666666666, 0.283071, 0.849213, 0.1223776112924619, 0.1223776112924619, 3, 1, 0.743817, 2.231451, 0.6787766238980854, 0.6787766238980855, 0.210553022007, 1.894977198063, 0.08306706183380948, 0.08306706183380949)
for ($c = 1; $c -le $row9.Length; $c++) {
    $ws.Cells.Item(9, $c).Value = $row9[$c-1]
}
